$p = $ppt.ActivePresentation
Write-Host ($p.Fonts | Get-Member)
